# Populate row 2 of the scoping document worksheet with NetSuite system data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NetSuite"
$ws.Range("B2").Value = "This system is primarily used for financial management, accounting, and reporting. It supports journal entries, trial balances, financial reporting, purchase order creation, and workflows for approvals."
$ws.Range("C2").Value = "The client utilizes this system to manage accounting processes, purchase order workflows, bill approvals, financial reporting, and integration with payment and budgeting tools."
$ws.Range("D2").Value = "The NetSuite system administration is managed by Blackforge Consulting, led by Steve Monti and Andrew Baker."
$ws.Range("E2").Value = "Access provisioning is initiated by an email sent by Megan Hodgson or another authorized individual, creating a ticket in Blackforge's system. The ticket requires approval from Megan Hodgson or John, with neither approving their own changes, before Blackforge executes the requested updates to roles or permissions and documents the changes in the system."
$ws.Range("F2").Value = "Upon termination, access is removed immediately. NetSuite access, integrated with Azure for Single Sign-On, is manually revoked, and a screenshot of the removal is documented and shared with HR via email."
$ws.Range("G2").Value = "Access is configured using a role-based model, where custom roles are defined and permissions are aligned with specific functions. Global permissions are disabled, and roles were tested during user acceptance testing as part of the system rollout."
$ws.Range("H2").Value = "Yes; System administrators can modify roles in NetSuite, including making changes to existing roles and testing these modifications as part of the user acceptance testing process."
$ws.Range("I2").Value = "Yes; management performs an annual review of all roles and permissions."
$ws.Range("J2").Value = "Yes; Users with privileged access include Steve Monti, Andrew Baker, and Justin, who hold the Administrator role and can make system changes."
$ws.Range("K2").Value = "Yes; there are system roles tied to individual user accounts, but no fully non-human accounts exist currently."
$ws.Range("L2").Value = "Credentials for shared and integration accounts are linked to individual user accounts and are not currently stored in a formal enterprise credential management solution."
$ws.Range("M2").Value = "N/A - This information was not discussed in the walkthrough meeting transcript."
$ws.Range("N2").Value = "Yes; Management performs periodic access reviews. Reviews were conducted during the initial implementation and go-live phase, with the intention to perform them quarterly moving forward."
$ws.Range("O2").Value = "The system maintains logs of administrative activities. Logs are retained indefinitely and include high-level user actions."
$ws.Range("P2").Value = "N/A - This information was not discussed in the walkthrough meeting transcript."
$ws.Range("Q2").Value = "The system uses Azure SSO for authentication in the production environment, while the sandbox environment relies on independent login functionality."
$ws.Range("R2").Value = "No; The client does not perform periodic reviews of the system's authentication configurations."
$ws.Range("S2").Value = "Management can perform configuration changes to adjust system settings, update workflows, modify integrations, and make limited code/script changes within the NetSuite system. These changes are managed through a ticketing system, with Blackforge handling development and deployment, and Jade's team performing testing and approvals where applicable."
$ws.Range("T2").Value = "Access to make changes is restricted to a small group of administrators. Steve Monti and Andrew Baker from Blackforge Consulting have full administrator access to NetSuite as primary administrators, while Justin from Whipfli currently has administrator access for managing integration setup. All change access is granted through role-based permissions in the system."
$ws.Range("U2").Value = "Management maintains separate Sandbox and Release Preview environments for implementing and testing changes."
$ws.Range("V2").Value = "The change management process varies by change type. Workflow, configuration, and script changes begin with ticket creation in Blackforge's system, followed by development in the sandbox environment, testing by Jade's team, approval by Jade's team, and deployment to production by Blackforge. Integration changes involve email requests to the integration owner, development by the owner or through NetSuite adjustments, testing in the live environment, informal approval via email, and direct deployment. Bundle audit trail changes and patches are managed by NetSuite, with Blackforge handling testing, approval, and deployment as needed. All changes follow distinct processes based on their type and risk level."
$ws.Range("W2").Value = "Yes, the vendor pushes updates, patches, and bug fixes. Semiannual system upgrades are tested in a release preview environment before deployment to Production, with evidence of testing documented. Patches and bug fixes are assessed for relevance and handled on a case-by-case basis."
$ws.Range("X2").Value = "No, there is no inherent system functionality preventing users from both developing and deploying changes."
$ws.Range("Y2").Value = "N/A - This information was not discussed in the walkthrough meeting transcript."
$ws.Range("Z2").Value = "Yes, the system has automated jobs and interfaces with Bill.com, Concur, and Adaptive Planning. These jobs perform data integration and transfer functions between systems."
$ws.Range("AA2").Value = "The following tools are used: NetSuite and Adaptive. Each tool is responsible for managing integration jobs and monitoring failures, with notifications sent for task success or failure."
$ws.Range("AB2").Value = "Job failures are handled through a defined process: detection occurs via Adaptive's run log, notifications are sent for successful or failed runs, and resolution involves identifying the issue, creating missing mappings, and manually re-running the job."
$ws.Range("AC2").Value = "The system utilizes NetSuite for data storage, which is vendor managed."
$ws.Range("AD2").Value = "N/A - This information was not discussed in the walkthrough meeting transcript."
$ws.Range("AE2").Value = "N/A - All backup types and strategies are managed by the vendor as part of their SaaS service."
$ws.Range("AF2").Value = "N/A - Backup failure resolution is managed by the vendor as part of their SaaS service."
$ws.Range("AG2").Value = "No, management does not perform regular SOC report reviews."
